# Auto-generated Excel COM-interop edit script.
# Updates the cryptos list (prices + 1h volume %) per the commit diff,
# including a name/link swap for rows 38-39 and a replaced coin in row 51.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cells whose new value reads as a plain number (e.g. "315.47") ---
# --- must be forced back to text so they stay inline/shared strings ---
# --- (matching the workbook's original plain-text storage) instead  ---
# --- of being auto-converted to a numeric cell by Excel's input parser. ---
$forcedTextValues = @{
    "D5" = "315.47"
    "D6" = "97.52"
    "D7" = "0.576"
    "D9" = "0.536"
    "D10" = "35.68"
    "D11" = "0.0813"
    "D12" = "7.51"
    "D16" = "15.26"
    "D17" = "0.846"
    "D19" = "6.83"
    "D20" = "12.69"
    "D22" = "69.52"
    "D23" = "251.58"
    "D26" = "27.37"
    "D29" = "41.05"
    "D30" = "10.30"
    "D31" = "5.86"
    "D32" = "157.22"
    "D35" = "0.0806"
    "D37" = "18.79"
    "D38" = "2.52"
    "D39" = "0.113"
    "D41" = "23.08"
    "D42" = "3.97"
    "D43" = "0.0305"
    "D47" = "8.98"
    "D49" = "83.31"
    "D51" = "74.37"
}
foreach ($addr in $forcedTextValues.Keys) {
    $ws.Range($addr).NumberFormat = "@"
}
foreach ($addr in $forcedTextValues.Keys) {
    $ws.Range($addr).Value = $forcedTextValues[$addr]
}
foreach ($addr in $forcedTextValues.Keys) {
    $ws.Range($addr).ClearFormats()
}

# --- Remaining cells (non-numeric-looking text: links, names, percents) ---
$ws.Range("D2").Value = "43.133.18"
$ws.Range("E2").Value = "  +1.45%  "
$ws.Range("D3").Value = "2.591.89"
$ws.Range("E3").Value = "  +3.08%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("E5").Value = "  -0.75%  "
$ws.Range("E6").Value = "  +3.48%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  +1.38%  "
$ws.Range("E10").Value = "  +0.14%  "
$ws.Range("E11").Value = "  +0.10%  "
$ws.Range("E12").Value = "  -0.18%  "
$ws.Range("D13").Value = "2.992.88"
$ws.Range("E13").Value = "  +3.13%  "
$ws.Range("E14").Value = "  -0.62%  "
$ws.Range("D15").Value = "2.588.53"
$ws.Range("E15").Value = "  +3.47%  "
$ws.Range("E16").Value = "  +0.34%  "
$ws.Range("E17").Value = "  -0.02%  "
$ws.Range("D18").Value = "43.261.41"
$ws.Range("E18").Value = "  +1.49%  "
$ws.Range("E19").Value = "  +2.72%  "
$ws.Range("E20").Value = "  -1.52%  "
$ws.Range("D21").Value = "0.0₃0967"
$ws.Range("E21").Value = "  +1.08%  "
$ws.Range("E22").Value = "  +0.35%  "
$ws.Range("E23").Value = "  +0.35%  "
$ws.Range("E24").Value = "  +0.11%  "
$ws.Range("E25").Value = "  +3.17%  "
$ws.Range("E26").Value = "  +2.68%  "
$ws.Range("E28").Value = "  -0.97%  "
$ws.Range("E29").Value = "  +0.37%  "
$ws.Range("E30").Value = "  +0.78%  "
$ws.Range("E31").Value = "  -0.78%  "
$ws.Range("E32").Value = "  +0.15%  "
$ws.Range("E33").Value = "  +5.82%  "
$ws.Range("E34").Value = "  +2.14%  "
$ws.Range("E35").Value = "  +3.60%  "
$ws.Range("E36").Value = "  +3.21%  "
$ws.Range("E37").Value = "  -1.40%  "
$ws.Range("B38").Value = "ApeXProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("E38").Value = "  +10.64%  "
$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("E39").Value = "  +1.93%  "
$ws.Range("E40").Value = "  +0.59%  "
$ws.Range("E41").Value = "  -1.33%  "
$ws.Range("E43").Value = "  +0.94%  "
$ws.Range("E44").Value = "  -0.14%  "
$ws.Range("D45").Value = "2.015.08"
$ws.Range("E45").Value = "  +0.13%  "
$ws.Range("E46").Value = "  -2.06%  "
$ws.Range("E47").Value = "  +0.97%  "
$ws.Range("D48").Value = "2.840.67"
$ws.Range("E48").Value = "  +3.02%  "
$ws.Range("E49").Value = "  -1.90%  "
$ws.Range("E50").Value = "  +4.45%  "
$ws.Range("B51").Value = "ordi"
$ws.Range("C51").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("E51").Value = "  -0.19%  "
